$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update the "resumen" fields (Valor Mora / Cant. Trabajadores / Cant. Periodos) ---
$ws.Range("E11").Value = 217880
$ws.Range("C13").Value = 2
$ws.Range("F13").Value = 4

# --- 2. Apply the "last row" (bottom-border) formatting from row 21 onto row 19, ---
#        since after deleting the two trailing worker rows, row 19 becomes the last
#        row of the table and must carry the thicker bottom border style.
$ws.Range("B21:J21").Copy()
$ws.Range("B19:J19").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

# --- 3. Rewrite the worker/period rows 16-19 with the refreshed data set ---
$ws.Range("E16").Value = "2503"
$ws.Range("G16").Value = 1423500

$ws.Range("E17").Value = "2506"
$ws.Range("G17").Value = 1423500

$ws.Range("C18").Value = "1143380075"
$ws.Range("D18").Value = "JHAXLYN NATALIA NARCISA ROCERO HERNANDEZ"
$ws.Range("E18").Value = "2507"
$ws.Range("F18").Value = 56940
$ws.Range("G18").Value = 1423500

$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "73009373"
$ws.Range("D19").Value = "XAVIER ELIAS TORRES MIRANDA"
$ws.Range("E19").Value = "2508"
$ws.Range("F19").Value = 56940
$ws.Range("G19").Value = 1423500

# --- 4. Remove the two now-obsolete worker rows (old rows 20 & 21); this shifts the ---
#        signature block (rows 26/27) up to rows 24/25 and keeps all merges in sync.
$ws.Range("B20:B21").EntireRow.Delete()
